# Handles float input without breaking stuff
# Fix the marksheet report: populate the "No."/"Marking"/"Total" summary row
# values with the real (non-zero) results, fill in the "Student Ans" column
# for each answered question (with correct/incorrect styling), and remove
# the two duplicate "Student Ans / Correct Ans" blocks (columns D:E below
# row 18, and the whole G:H block) that the buggy report generator used to
# emit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Summary block (rows 10-12): give the row-label cells (A10:A12) the
#    same "mtitleStyle" used by the header row above them, and update the
#    numeric results.
# ---------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

# No. row
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 28

# Marking row
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Total row
$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "86/112"

# ---------------------------------------------------------------------
# 2. Remove the duplicate "Student Ans / Correct Ans" report blocks.
# ---------------------------------------------------------------------
# Whole third block (columns G:H, rows 15-40)
$ws.Range("G15:H40").Clear()
# Second block only keeps its header + first three question rows (16-18);
# everything from row 19 down is removed entirely.
$ws.Range("D19:E40").Clear()

# ---------------------------------------------------------------------
# 3. Fill in the "Student Ans" column (A) for the first block, and the
#    surviving rows of the second block (D16:D18), copying the
#    correct/incorrect/normal style from existing cells that already use
#    it so no new duplicate styles are created.
# ---------------------------------------------------------------------
$ws.Range("B10").Copy()   # correctStyle source
$correctTargets = "A16,A17,A18,A19,A21,A22,A23,A24,A27,A28,A29,A30,A31,A32,A33,A35,A37,A38,A39,D16,D17,D18"
$ws.Range($correctTargets).PasteSpecial(-4122)

$ws.Range("C10").Copy()   # incorrectStyle source
$ws.Range("A34,A36").PasteSpecial(-4122)

# Values for column A (first block)
$ws.Range("A16").Value = "Option A"
$ws.Range("A17").Value = "Option D"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A23").Value = "Option D"
$ws.Range("A24").Value = "Option A"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A34").Value = "Option A"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option D"
$ws.Range("A37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"

# Values for column D (surviving rows of the second block)
$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"
